$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-24 Monday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-06-25 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("14×39=546", $true, $true, $false, $false, $false, $true, 1, $false, "59×45=2655", 2) | Out-Null
$d.Content.Find.Execute("85×56=4760", $true, $true, $false, $false, $false, $true, 1, $false, "67×95=6365", 2) | Out-Null
$d.Content.Find.Execute("21×74=1554", $true, $true, $false, $false, $false, $true, 1, $false, "96×80=7680", 2) | Out-Null
$d.Content.Find.Execute("99×92=9108", $true, $true, $false, $false, $false, $true, 1, $false, "27×26=702", 2) | Out-Null
$d.Content.Find.Execute("44×81=3564", $true, $true, $false, $false, $false, $true, 1, $false, "99×38=3762", 2) | Out-Null
$d.Content.Find.Execute("40×60=2400", $true, $true, $false, $false, $false, $true, 1, $false, "99×46=4554", 2) | Out-Null
$d.Content.Find.Execute("59×58=3422", $true, $true, $false, $false, $false, $true, 1, $false, "62×82=5084", 2) | Out-Null
$d.Content.Find.Execute("29×28=812", $true, $true, $false, $false, $false, $true, 1, $false, "61×56=3416", 2) | Out-Null
$d.Content.Find.Execute("87×44=3828", $true, $true, $false, $false, $false, $true, 1, $false, "55×79=4345", 2) | Out-Null
$d.Content.Find.Execute("84×43=3612", $true, $true, $false, $false, $false, $true, 1, $false, "92×95=8740", 2) | Out-Null
$d.Content.Find.Execute("47×42=1974", $true, $true, $false, $false, $false, $true, 1, $false, "94×89=8366", 2) | Out-Null
$d.Content.Find.Execute("77×88=6776", $true, $true, $false, $false, $false, $true, 1, $false, "28×62=1736", 2) | Out-Null
$d.Content.Find.Execute("43×24=1032", $true, $true, $false, $false, $false, $true, 1, $false, "20×88=1760", 2) | Out-Null
$d.Content.Find.Execute("19×22=418", $true, $true, $false, $false, $false, $true, 1, $false, "86×75=6450", 2) | Out-Null
$d.Content.Find.Execute("66×82=5412", $true, $true, $false, $false, $false, $true, 1, $false, "60×99=5940", 2) | Out-Null
$d.Content.Find.Execute("21×69=1449", $true, $true, $false, $false, $false, $true, 1, $false, "38×99=3762", 2) | Out-Null
$d.Content.Find.Execute("86×70=6020", $true, $true, $false, $false, $false, $true, 1, $false, "18×83=1494", 2) | Out-Null
$d.Content.Find.Execute("15×82=1230", $true, $true, $false, $false, $false, $true, 1, $false, "60×98=5880", 2) | Out-Null
$d.Content.Find.Execute("15×81=1215", $true, $true, $false, $false, $false, $true, 1, $false, "16×97=1552", 2) | Out-Null
$d.Content.Find.Execute("86×94=8084", $true, $true, $false, $false, $false, $true, 1, $false, "26×14=364", 2) | Out-Null
$d.Content.Find.Execute("11×85=935", $true, $true, $false, $false, $false, $true, 1, $false, "18×30=540", 2) | Out-Null
$d.Content.Find.Execute("35×13=455", $true, $true, $false, $false, $false, $true, 1, $false, "16×30=480", 2) | Out-Null
$d.Content.Find.Execute("66×42=2772", $true, $true, $false, $false, $false, $true, 1, $false, "73×23=1679", 2) | Out-Null
$d.Content.Find.Execute("66×16=1056", $true, $true, $false, $false, $false, $true, 1, $false, "73×60=4380", 2) | Out-Null
$d.Content.Find.Execute("56×91=5096", $true, $true, $false, $false, $false, $true, 1, $false, "27×41=1107", 2) | Out-Null
